$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = 5.569900000000001
    12 = 5.031799999999999
    15 = 4.625499999999997
    27 = 6.552199999999999
    28 = 6.417700000000004
    31 = 5.813900000000001
    32 = 6.583099999999997
    36 = 9.191700000000003
    38 = 5.1732
    46 = 7.040600000000003
    54 = 4.677600000000002
    55 = 5.133699999999999
    56 = 4.865199999999999
    67 = 5.562599999999995
    69 = 5.354699999999997
    72 = 5.166600000000003
    73 = 8.473099999999997
    83 = 5.507699999999997
    86 = 4.833900000000002
    91 = 5.416600000000002
    93 = 5.753400000000003
    99 = 4.576699999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
